$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expenses")

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-03-24"
$ws.Range("A3").NumberFormat = "General"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 350
$ws.Range("C3").Value = "Shopping at"
$ws.Range("D3").Value = "kzon"
